$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new paragraph "- Acesso remoto ao sistema" right
#    before the "- Controle de pagamentos com multiplos meios(...)"
#    paragraph, re-using that paragraph's bold / sz28 formatting.
# -----------------------------------------------------------------
$pControle = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- Controle de pagamentos com m*") {
        $pControle = $p
    }
}
$pControle.Range.InsertParagraphBefore() | Out-Null

# The paragraph collection shifted, so re-locate the target paragraph
# and grab the freshly created (empty) one right before it.
$pAcesso = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- Controle de pagamentos com m*") {
        $pAcesso = $p.Previous()
    }
}
$pAcesso.Range.Text = "- Acesso remoto ao sistema"

# -----------------------------------------------------------------
# 2) Insert a new paragraph "- Emissao do ticket do estacionamento"
#    right after the "Controle de pagamentos..." paragraph.
# -----------------------------------------------------------------
$pControle2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- Controle de pagamentos com m*") {
        $pControle2 = $p
    }
}
$pControle2.Range.InsertParagraphAfter() | Out-Null

$pTicket = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "- Controle de pagamentos com m*") {
        $pTicket = $p.Next()
    }
}
$pTicket.Range.Text = "- Emissão do ticket do estacionamento"

# Move the (hidden) _GoBack bookmark from the end of the
# "Cadastramento de politica..." paragraph onto the end of this new
# "Emissao do ticket..." paragraph (collapsed, right before the
# paragraph mark - matching where it originally sat).
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

$goBackRange = $pTicket.Range.Duplicate
$goBackRange.MoveEnd(1, -1) | Out-Null
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# -----------------------------------------------------------------
# 3) Drop the stray <w:lastRenderedPageBreak/> that sits in front of
#    "- Brainstorm:" — rewriting the run's text regenerates the run
#    without the rendering-only breadcrumb.
# -----------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Brainstorm*") {
        $p.Range.Text = "- Brainstorm:"
    }
}
